$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records need to be inserted above the existing
# "Femacal de La Calera - Papaya" history (row 64 onwards), pushing the
# rest of the table down by two rows.
$ws.Rows("64:65").Insert()

# Seed the two freshly-inserted rows with the same layout/formatting as
# the (now shifted) rows immediately below them, then overwrite the
# values that actually differ for this week's entries.
$ws.Range("A66:T66").Copy()
$ws.Range("A64").PasteSpecial()
$ws.Range("A67:T67").Copy()
$ws.Range("A65").PasteSpecial()

# Row 64: Fecha 2023-07-20 (serial 45127), Volumen 50 (Primera)
$ws.Range("D64").Value = 45127
$ws.Range("M64").Value = 50

# Row 65: Fecha 2023-07-20 (serial 45127), Volumen 36 (Segunda)
$ws.Range("D65").Value = 45127
$ws.Range("M65").Value = 36
